$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

# Columns A (Date) and D (Week) hold text that looks like a date / a plain
# number ("2023-06-27" and "26"). Force them to be stored as literal text
# (matching the rest of the column) instead of being auto-converted to a
# date serial / number by temporarily marking the cell as Text before the
# assignment, then clearing the number format again so no extra formatting
# is left behind on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-27"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "15:24:33"
$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 122846
$ws.Cells.Item($row, 6).Value = 134293
$ws.Cells.Item($row, 7).Value = 163473
$ws.Cells.Item($row, 8).Value = 133815
$ws.Cells.Item($row, 9).Value = 177259
$ws.Cells.Item($row, 10).Value = 114970
$ws.Cells.Item($row, 11).Value = 203716
$ws.Cells.Item($row, 12).Value = 226268
$ws.Cells.Item($row, 13).Value = 176236
$ws.Cells.Item($row, 14).Value = 104336
$ws.Cells.Item($row, 15).Value = 39642
$ws.Cells.Item($row, 16).Value = 33755
$ws.Cells.Item($row, 17).Value = 52242
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35751
$ws.Cells.Item($row, 20).Value = -1
